# Add a new row of data (row 29) to Sheet1, mirroring the formatting of row 28.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row index
$newRow = 29

# A29: numeric value 27, same style/formatting as A28 (bold/border/centered)
$ws.Cells.Item($newRow - 1, 1).Copy($ws.Cells.Item($newRow, 1))
$ws.Cells.Item($newRow, 1).Value = 27

# B29: inline string ".7."
$ws.Cells.Item($newRow, 2).Value = ".7."

# C29: numeric value 187
$ws.Cells.Item($newRow, 3).Value = 187

$wb.Save()
